# Auto-generated edit script applying the crypto price/volume update diff.
# Updates coin names/links (text) and price/volume values (kept as literal text
# to match the source workbook's inlineStr storage) for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, [string]$cellRef, [string]$value) {
    # Force the cell to Text format first so Excel stores the numeric-looking
    # string (price / percentage) as literal text instead of converting it
    # to a Number/Percentage value.
    $range = $sheet.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws 'D2' '303.98'
Set-TextValue $ws 'E2' '3.25%'
Set-TextValue $ws 'D3' '43.64'
Set-TextValue $ws 'E3' '8.64%'
Set-TextValue $ws 'D4' '5.076'
Set-TextValue $ws 'E4' '1.25%'
Set-TextValue $ws 'D5' '0.07675'
Set-TextValue $ws 'E5' '4.35%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws 'D6' '1.610'
Set-TextValue $ws 'E6' '3.76%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D7' '1.004'
Set-TextValue $ws 'E7' '8.63%'
$ws.Range('B8').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C8').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D8' '0.1253'
Set-TextValue $ws 'E8' '7.50%'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D9' '0.1855'
Set-TextValue $ws 'E9' '3.42%'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D10' '0.09128'
Set-TextValue $ws 'E10' '4.81%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D11' '0.04166'
Set-TextValue $ws 'E11' '-2.56%'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D12' '0.1046'
Set-TextValue $ws 'E12' '-0.51%'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D13' '0.001268'
Set-TextValue $ws 'E13' '-0.64%'
$ws.Range('B14').Value = 'TigerCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D14' '0.005745'
Set-TextValue $ws 'E14' '-4.30%'
$ws.Range('B15').Value = 'UpBots'
$ws.Range('C15').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws 'D15' '0.007430'
Set-TextValue $ws 'E15' '1,897.12%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D16' '3.329'
Set-TextValue $ws 'E16' '-0.43%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D17' '4.411'
Set-TextValue $ws 'E17' '2.26%'
Set-TextValue $ws 'D19' '0.3354'
Set-TextValue $ws 'E19' '1.85%'
Set-TextValue $ws 'D20' '8.449'
Set-TextValue $ws 'E20' '6.70%'
Set-TextValue $ws 'D21' '0.1399'
Set-TextValue $ws 'E21' '1.09%'
Set-TextValue $ws 'E22' '-3.06%'
Set-TextValue $ws 'D23' '0.04166'
Set-TextValue $ws 'E23' '5.71%'
Set-TextValue $ws 'E24' '1.11%'
Set-TextValue $ws 'D25' '0.004495'
Set-TextValue $ws 'E25' '18.43%'
Set-TextValue $ws 'D26' '0.0001347'
Set-TextValue $ws 'E26' '9.56%'
Set-TextValue $ws 'D38' '0.02454'
Set-TextValue $ws 'E38' '4.64%'
Set-TextValue $ws 'D39' '0.05294'
Set-TextValue $ws 'D40' '0.005958'
Set-TextValue $ws 'E40' '2.76%'
Set-TextValue $ws 'D41' '0.007644'
Set-TextValue $ws 'E41' '-1.34%'
Set-TextValue $ws 'D42' '0.1347'
Set-TextValue $ws 'E42' '4.02%'
Set-TextValue $ws 'D43' '0.007357'
Set-TextValue $ws 'E43' '-0.22%'
Set-TextValue $ws 'D44' '0.007551'
Set-TextValue $ws 'E44' '8.28%'
Set-TextValue $ws 'D45' '0.3026'
Set-TextValue $ws 'E45' '3.51%'
Set-TextValue $ws 'D46' '0.00006712'
Set-TextValue $ws 'E46' '7.97%'
Set-TextValue $ws 'E47' '-0.17%'
Set-TextValue $ws 'D48' '0.04173'
Set-TextValue $ws 'E48' '-10.01%'
Set-TextValue $ws 'E49' '0.05%'
Set-TextValue $ws 'E50' '-0.17%'
Set-TextValue $ws 'E51' '-0.17%'
